$d = $word.ActiveDocument

# --- 1. Fix the "Dockefile" -> "Dockerfile" typo -------------------------
# Insert the missing "r" right after "Docke" (keeps the existing
# w:proofErr spellStart/spellEnd wrapper intact around the word).
$rng = $d.Content
$found = $rng.Find.Execute("Docke", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPt = $d.Range($rng.End, $rng.End)
    $insertPt.InsertAfter("r")
}

# --- 2. Move the "_GoBack" bookmark ---------------------------------------
# It currently sits after "...hello.py inside folder CONTAINER"; it should
# instead sit inside "...to be able | to build this image" (right after
# "to be able").
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$rng2 = $d.Content
$found2 = $rng2.Find.Execute(" to be able", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $bmPoint = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}

# --- 3. Mark every inline picture's run as "do not spell/grammar check" --
# (adds <w:noProof/> to the run properties of each drawing run, matching
# what Word stamps on pasted/inserted pictures.)
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shapeRange = $d.InlineShapes($i).Range
    $shapeRange.NoProofing = -1
}
